$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Clear-Cell($row, $col) {
    $ws.Cells.Item($row, $col).Value = ""
}

Set-TextCell 4 8 "./instructions_videos/block_4_text.mp4"  # H4
$ws.Cells.Item(4, 9).Value = 4  # I4
Set-TextCell 5 5 "2"  # E5
Set-TextCell 5 7 "inverse"  # G5
Set-TextCell 5 8 "../stimuli/exp_videos/VR/2.mp4"  # H5
$ws.Cells.Item(5, 9).Value = 4  # I5
Set-TextCell 6 8 "./instructions_videos/post_stimulus_self_report.mp4"  # H6
$ws.Cells.Item(6, 9).Value = 4  # I6
Set-TextCell 6 12 "post_stimulus_self_report"  # L6
Set-TextCell 7 8 "./instructions_videos/mareo.mp4"  # H7
$ws.Cells.Item(7, 9).Value = 4  # I7
Set-TextCell 7 12 "motion_sickness"  # L7
Set-TextCell 8 8 "./instructions_videos/block_4_text_reminder.mp4"  # H8
$ws.Cells.Item(8, 9).Value = 4  # I8
Set-TextCell 8 12 "audio_instruction"  # L8
$ws.Cells.Item(9, 4).Value = 2  # D9
Set-TextCell 9 5 "12"  # E9
Set-TextCell 9 6 "valence"  # F9
Set-TextCell 9 7 "inverse"  # G9
Set-TextCell 9 8 "../stimuli/exp_videos/VR/12.mp4"  # H9
$ws.Cells.Item(9, 9).Value = 4  # I9
Set-TextCell 9 12 "video"  # L9
Set-TextCell 10 8 "./instructions_videos/post_stimulus_self_report.mp4"  # H10
$ws.Cells.Item(10, 9).Value = 4  # I10
Set-TextCell 10 12 "post_stimulus_self_report"  # L10
Clear-Cell 11 4  # D11: was 2
Clear-Cell 11 5  # E11: was '6'
Clear-Cell 11 6  # F11: was 'valence'
Clear-Cell 11 7  # G11: was 'direct'
Set-TextCell 11 8 "./instructions_videos/mareo.mp4"  # H11
$ws.Cells.Item(11, 9).Value = 4  # I11
Set-TextCell 11 12 "motion_sickness"  # L11
Set-TextCell 12 8 "./instructions_videos/block_4_text_reminder.mp4"  # H12
$ws.Cells.Item(12, 9).Value = 4  # I12
Set-TextCell 12 12 "audio_instruction"  # L12
$ws.Cells.Item(13, 4).Value = 3  # D13
Set-TextCell 13 5 "3"  # E13
Set-TextCell 13 6 "valence"  # F13
Set-TextCell 13 7 "inverse"  # G13
Set-TextCell 13 8 "../stimuli/exp_videos/VR/3.mp4"  # H13
$ws.Cells.Item(13, 9).Value = 4  # I13
Set-TextCell 13 12 "video"  # L13
Set-TextCell 14 8 "./instructions_videos/post_stimulus_self_report.mp4"  # H14
$ws.Cells.Item(14, 9).Value = 4  # I14
Set-TextCell 14 12 "post_stimulus_self_report"  # L14
Set-TextCell 15 8 "./instructions_videos/mareo.mp4"  # H15
$ws.Cells.Item(15, 9).Value = 4  # I15
Set-TextCell 15 12 "motion_sickness"  # L15
Set-TextCell 16 8 "./instructions_videos/luminance_instructions_inverse.mp4"  # H16
Clear-Cell 16 9  # I16: was 3
Set-TextCell 16 12 "luminance_instructions"  # L16
$ws.Cells.Item(17, 4).Value = 4  # D17
Clear-Cell 17 5  # E17: was '14'
Set-TextCell 17 6 "luminance"  # F17
Set-TextCell 17 7 "inverse"  # G17
Set-TextCell 17 8 "../stimuli/exp_videos/VR/green_intensity_video_12.mp4"  # H17
Clear-Cell 17 9  # I17: was 3
Set-TextCell 17 12 "luminance"  # L17
Set-TextCell 18 8 "./instructions_videos/confidence_luminance_practice_instructions_text.mp4"  # H18
Clear-Cell 18 9  # I18: was 3
Set-TextCell 18 12 "confidence_luminance_instructions"  # L18
Set-TextCell 19 8 "./instructions_videos/block_1_text.mp4"  # H19
$ws.Cells.Item(19, 9).Value = 1  # I19
$ws.Cells.Item(19, 10).Value = 1  # J19
Set-TextCell 19 12 "audio_instruction"  # L19
$ws.Cells.Item(20, 4).Value = 5  # D20
Set-TextCell 20 5 "4"  # E20
Set-TextCell 20 6 "arousal"  # F20
Set-TextCell 20 7 "direct"  # G20
Set-TextCell 20 8 "../stimuli/exp_videos/VR/4.mp4"  # H20
$ws.Cells.Item(20, 9).Value = 1  # I20
$ws.Cells.Item(20, 10).Value = 1  # J20
Set-TextCell 20 12 "video"  # L20
Set-TextCell 21 8 "./instructions_videos/post_stimulus_verbal_report.mp4"  # H21
$ws.Cells.Item(21, 9).Value = 1  # I21
$ws.Cells.Item(21, 10).Value = 1  # J21
Set-TextCell 21 12 "instruction_post_stimulus_verbal_report"  # L21
Set-TextCell 22 8 "./videos_fixation/countdown_bar.mp4"  # H22
$ws.Cells.Item(22, 9).Value = 1  # I22
$ws.Cells.Item(22, 10).Value = 1  # J22
Set-TextCell 22 12 "verbal_report"  # L22
Clear-Cell 23 4  # D23: was 4
Clear-Cell 23 6  # F23: was 'luminance'
Clear-Cell 23 7  # G23: was 'direct'
Set-TextCell 23 8 "./instructions_videos/confidence_verbal_report_text.mp4"  # H23
$ws.Cells.Item(23, 9).Value = 1  # I23
$ws.Cells.Item(23, 10).Value = 1  # J23
Set-TextCell 23 12 "confidence_verbal_report"  # L23
Set-TextCell 24 8 "./instructions_videos/mareo.mp4"  # H24
$ws.Cells.Item(24, 9).Value = 1  # I24
Set-TextCell 24 12 "motion_sickness"  # L24
Clear-Cell 25 4  # D25: was 5
Clear-Cell 25 5  # E25: was '1'
Clear-Cell 25 6  # F25: was 'arousal'
Clear-Cell 25 7  # G25: was 'inverse '
Set-TextCell 25 8 "./instructions_videos/block_1_text_reminder.mp4"  # H25
$ws.Cells.Item(25, 9).Value = 1  # I25
Set-TextCell 25 12 "audio_instruction"  # L25
$ws.Cells.Item(26, 4).Value = 6  # D26
Set-TextCell 26 5 "9"  # E26
Set-TextCell 26 6 "arousal"  # F26
Set-TextCell 26 7 "direct"  # G26
Set-TextCell 26 8 "../stimuli/exp_videos/VR/9.mp4"  # H26
$ws.Cells.Item(26, 9).Value = 1  # I26
Set-TextCell 26 12 "video"  # L26
Set-TextCell 27 8 "./instructions_videos/post_stimulus_verbal_report.mp4"  # H27
$ws.Cells.Item(27, 9).Value = 1  # I27
Set-TextCell 27 12 "instruction_post_stimulus_verbal_report"  # L27
Set-TextCell 28 8 "./videos_fixation/countdown_bar.mp4"  # H28
$ws.Cells.Item(28, 9).Value = 1  # I28
Set-TextCell 28 12 "verbal_report"  # L28
Clear-Cell 29 4  # D29: was 6
Clear-Cell 29 5  # E29: was '10'
Clear-Cell 29 6  # F29: was 'arousal'
Clear-Cell 29 7  # G29: was 'inverse'
Set-TextCell 29 8 "./instructions_videos/confidence_verbal_report_text.mp4"  # H29
$ws.Cells.Item(29, 9).Value = 1  # I29
Set-TextCell 29 12 "confidence_verbal_report"  # L29
Set-TextCell 30 8 "./instructions_videos/mareo.mp4"  # H30
$ws.Cells.Item(30, 9).Value = 1  # I30
Set-TextCell 30 12 "motion_sickness"  # L30
Set-TextCell 31 8 "./instructions_videos/block_1_text_reminder.mp4"  # H31
$ws.Cells.Item(31, 9).Value = 1  # I31
Set-TextCell 31 12 "audio_instruction"  # L31
$ws.Cells.Item(32, 4).Value = 7  # D32
Set-TextCell 32 5 "7"  # E32
Set-TextCell 32 6 "arousal"  # F32
Set-TextCell 32 7 "direct"  # G32
Set-TextCell 32 8 "../stimuli/exp_videos/VR/7.mp4"  # H32
$ws.Cells.Item(32, 9).Value = 1  # I32
Set-TextCell 32 12 "video"  # L32
Clear-Cell 33 4  # D33: was 7
Clear-Cell 33 5  # E33: was '5'
Clear-Cell 33 6  # F33: was 'arousal'
Clear-Cell 33 7  # G33: was 'inverse'
Set-TextCell 33 8 "./instructions_videos/post_stimulus_verbal_report.mp4"  # H33
$ws.Cells.Item(33, 9).Value = 1  # I33
Set-TextCell 33 12 "instruction_post_stimulus_verbal_report"  # L33
Set-TextCell 34 8 "./videos_fixation/countdown_bar.mp4"  # H34
$ws.Cells.Item(34, 9).Value = 1  # I34
Set-TextCell 34 12 "verbal_report"  # L34
Set-TextCell 35 8 "./instructions_videos/confidence_verbal_report_text.mp4"  # H35
$ws.Cells.Item(35, 9).Value = 1  # I35
Set-TextCell 35 12 "confidence_verbal_report"  # L35
Set-TextCell 36 8 "./instructions_videos/mareo.mp4"  # H36
$ws.Cells.Item(36, 9).Value = 1  # I36
Set-TextCell 36 12 "motion_sickness"  # L36
Clear-Cell 37 4  # D37: was 8
Clear-Cell 37 5  # E37: was '11'
Clear-Cell 37 6  # F37: was 'arousal'
Clear-Cell 37 7  # G37: was 'inverse'
Set-TextCell 37 8 "./instructions_videos/block_1_text_reminder.mp4"  # H37
$ws.Cells.Item(37, 9).Value = 1  # I37
Set-TextCell 37 12 "audio_instruction"  # L37
$ws.Cells.Item(38, 4).Value = 8  # D38
Set-TextCell 38 5 "8"  # E38
Set-TextCell 38 6 "arousal"  # F38
Set-TextCell 38 7 "direct"  # G38
Set-TextCell 38 8 "../stimuli/exp_videos/VR/8.mp4"  # H38
$ws.Cells.Item(38, 9).Value = 1  # I38
Set-TextCell 38 12 "video"  # L38
Set-TextCell 39 8 "./instructions_videos/post_stimulus_verbal_report.mp4"  # H39
$ws.Cells.Item(39, 9).Value = 1  # I39
Set-TextCell 39 12 "instruction_post_stimulus_verbal_report"  # L39
Set-TextCell 40 8 "./videos_fixation/countdown_bar.mp4"  # H40
$ws.Cells.Item(40, 9).Value = 1  # I40
Set-TextCell 40 12 "verbal_report"  # L40
Clear-Cell 41 4  # D41: was 9
Clear-Cell 41 6  # F41: was 'luminance'
Clear-Cell 41 7  # G41: was 'inverse '
Set-TextCell 41 8 "./instructions_videos/confidence_verbal_report_text.mp4"  # H41
$ws.Cells.Item(41, 9).Value = 1  # I41
Set-TextCell 41 12 "confidence_verbal_report"  # L41
Set-TextCell 42 8 "./instructions_videos/mareo.mp4"  # H42
$ws.Cells.Item(42, 9).Value = 1  # I42
Set-TextCell 42 12 "motion_sickness"  # L42
Set-TextCell 43 1 "12"  # A43
Set-TextCell 43 2 "A_block1"  # B43
Set-TextCell 43 3 "VR"  # C43
Set-TextCell 43 8 "./instructions_videos/luminance_instructions_direct.mp4"  # H43
$ws.Cells.Item(43, 10).Value = 1  # J43
$ws.Cells.Item(43, 11).Value = 1  # K43
Set-TextCell 43 12 "luminance_instructions"  # L43
Set-TextCell 44 1 "12"  # A44
Set-TextCell 44 2 "A_block1"  # B44
Set-TextCell 44 3 "VR"  # C44
$ws.Cells.Item(44, 4).Value = 9  # D44
Set-TextCell 44 6 "luminance"  # F44
Set-TextCell 44 7 "direct"  # G44
Set-TextCell 44 8 "../stimuli/exp_videos/VR/green_intensity_video_3.mp4"  # H44
$ws.Cells.Item(44, 10).Value = 1  # J44
$ws.Cells.Item(44, 11).Value = 1  # K44
Set-TextCell 44 12 "luminance"  # L44
Set-TextCell 45 1 "12"  # A45
Set-TextCell 45 2 "A_block1"  # B45
Set-TextCell 45 3 "VR"  # C45
Set-TextCell 45 8 "./instructions_videos/confidence_luminance_practice_instructions_text.mp4"  # H45
$ws.Cells.Item(45, 10).Value = 1  # J45
$ws.Cells.Item(45, 11).Value = 1  # K45
Set-TextCell 45 12 "confidence_luminance_instructions"  # L45
Set-TextCell 46 1 "12"  # A46
Set-TextCell 46 2 "A_block1"  # B46
Set-TextCell 46 3 "VR"  # C46
Set-TextCell 46 8 "./instructions_videos/rest_suprablock_text.mp4"  # H46
$ws.Cells.Item(46, 10).Value = 1  # J46
$ws.Cells.Item(46, 11).Value = 1  # K46
Set-TextCell 46 12 "rest_suprablock"  # L46
